# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly-scraped counts, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Mapping of row -> new value for column F, shared between the two sheets
# (the "全部类型" sheet aggregates rows from "展览", "演出" and "本地生活",
# so the same exhibitions appear there at slightly shifted row numbers).

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetExhibition.Range("F2").Value = 10076
$sheetExhibition.Range("F5").Value = 20
$sheetExhibition.Range("F8").Value = 478
$sheetExhibition.Range("F9").Value = 757
$sheetExhibition.Range("F12").Value = 1040
$sheetExhibition.Range("F13").Value = 3130
$sheetExhibition.Range("F14").Value = 2347
$sheetExhibition.Range("F16").Value = 2064
$sheetExhibition.Range("F18").Value = 1942
$sheetExhibition.Range("F19").Value = 484
$sheetExhibition.Range("F21").Value = 548
$sheetExhibition.Range("F22").Value = 51
$sheetExhibition.Range("F23").Value = 234
$sheetExhibition.Range("F24").Value = 4
$sheetExhibition.Range("F25").Value = 10
$sheetExhibition.Range("F27").Value = 43
$sheetExhibition.Range("F30").Value = 357
$sheetExhibition.Range("F31").Value = 573
$sheetExhibition.Range("F32").Value = 45
$sheetExhibition.Range("F33").Value = 226
$sheetExhibition.Range("F34").Value = 1566
$sheetExhibition.Range("F35").Value = 22
$sheetExhibition.Range("F36").Value = 309
$sheetExhibition.Range("F37").Value = 1650
$sheetExhibition.Range("F38").Value = 104
$sheetExhibition.Range("F39").Value = 415
$sheetExhibition.Range("F41").Value = 433
$sheetExhibition.Range("F42").Value = 933

$sheetAll = $wb.Worksheets.Item("全部类型")
$sheetAll.Range("F2").Value = 10076
$sheetAll.Range("F6").Value = 20
$sheetAll.Range("F10").Value = 478
$sheetAll.Range("F11").Value = 757
$sheetAll.Range("F13").Value = 1040
$sheetAll.Range("F14").Value = 3130
$sheetAll.Range("F15").Value = 2347
$sheetAll.Range("F16").Value = 2064
$sheetAll.Range("F17").Value = 2064
$sheetAll.Range("F18").Value = 1942
$sheetAll.Range("F19").Value = 484
$sheetAll.Range("F21").Value = 548
$sheetAll.Range("F22").Value = 51
$sheetAll.Range("F23").Value = 234
$sheetAll.Range("F24").Value = 4
$sheetAll.Range("F25").Value = 10
$sheetAll.Range("F27").Value = 43
$sheetAll.Range("F30").Value = 357
$sheetAll.Range("F31").Value = 573
$sheetAll.Range("F35").Value = 45
$sheetAll.Range("F36").Value = 226
$sheetAll.Range("F37").Value = 1566
$sheetAll.Range("F38").Value = 22
$sheetAll.Range("F40").Value = 309
$sheetAll.Range("F41").Value = 1650
$sheetAll.Range("F42").Value = 104
$sheetAll.Range("F44").Value = 415
$sheetAll.Range("F46").Value = 433
$sheetAll.Range("F47").Value = 933
